$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "No" column (column A); everything else shifts one column left
# (Pekerjaan Utama -> A, Uraian Pekerjaan -> B, Durasi (Hari) -> C, Bobot % -> D,
# "Hari ke ->" label + the 1..10 day numbers -> E..O).
$ws.Columns.Item(1).Delete()

# Replace the "Hari ke ->" label and the old 1..10 day-number series with a
# continuous 1..90 day-number series starting at column E.
for ($i = 1; $i -le 90; $i++) {
    $ws.Cells.Item(1, 4 + $i).Value = $i
}

# Match the saved selection state.
$null = $ws.Range("J15").Select()
